$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.142.16"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "4.023.75"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.61"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.78"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.700"
$ws.Range("E7").Value = "  +12.07%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.748"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000326"
$ws.Range("E11").Value = "  -4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.81"
$ws.Range("E12").Value = "  +5.48%  "
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "4.671.48"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "4.028.80"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.11"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.59"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").Value = "72.153.12"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.30"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "98.07"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.21"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.28"
$ws.Range("E26").Value = "  -7.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.71"
$ws.Range("E27").Value = "  -5.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.86"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.67"
$ws.Range("E29").Value = "  +19.41%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.72"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("E31").Value = "  +6.88%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.132"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.46"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "680.18"
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "47.94"
$ws.Range("E35").Value = "  +16.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.64"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.450"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "0.0₃0825"
$ws.Range("E38").Value = "  -9.54%  "
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.39"
$ws.Range("E40").Value = "  -7.53%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.30"
$ws.Range("E42").Value = "  +6.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0493"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").Value = "THORChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.25"
$ws.Range("E45").Value = "  +11.42%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.150"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.03"
$ws.Range("E49").Value = "  -4.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000267"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.25"
$ws.Range("E51").Value = "  -2.89%  "
